$wb = $excel.ActiveWorkbook

# "Rules" is the only (active) worksheet in this workbook.
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently shows the rule-row label "R40" (pulled from the shared
# string table). Replace it with the text "1", keeping it as literal text
# (not a number) by using the classic leading-apostrophe text qualifier -
# this is how Excel stores a number-look-alike as a string.
$ws.Range("B11").Value = "'1"
